$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 921 (shifts existing rows 921-946 down to 922-947)
$ws.Rows.Item(921).Insert()

# Populate the newly inserted row 921 with the new weekly price record
$ws.Range("A921").Value = 8
$ws.Range("B921").Value = "Terminal La Palmera de La Serena"
$ws.Range("C921").Value = "Coquimbo"
$ws.Range("D921").Value = 45239
$ws.Range("D921").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E921").Value = 4
$ws.Range("F921").Value = 100112024
$ws.Range("G921").Value = "Choclo"
$ws.Range("H921").Value = "Dulce o Americano"
$ws.Range("I921").Value = "Primera"
$ws.Range("J921").Value = 400
$ws.Range("K921").Value = 28000
$ws.Range("L921").Value = 29000
$ws.Range("M921").Value = 28500
$ws.Range("N921").Value = "$/malla 70 unidades"
$ws.Range("O921").Value = "Región de Arica y Parinacota"
$ws.Range("P921").Value = 407
$ws.Range("Q921").Value = 70
$ws.Range("R921").Value = "Hortaliza"
